# Saldo_guide.xlsx update: refresh the "Dt. Referencia" (column G) from
# 2024-07-11 (serial 45484) to 2024-07-15 (serial 45488) for every data row,
# correct three mis-scaled balances (columns E / H), and rename the sheet
# (and its tab) to match the new extraction timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the reference date for all data rows (2-275) -------------------
$ws.Range("G2:G275").Value2 = 45488

# --- Correct the three balances that were off by a factor of ~1000 ---------
$ws.Range("E55").Value2  = 143371.44
$ws.Range("H55").Value2  = 143371.44

$ws.Range("E102").Value2 = 11726.16
$ws.Range("H102").Value2 = 11726.16

$ws.Range("E103").Value2 = 21904.13
$ws.Range("H103").Value2 = 21904.13

# --- Rename the sheet to match the new extraction run timestamp ------------
$ws.Name = "IClientBalance-20240715-094706-"
